# Auto-generated script to apply profit recalculation updates
# to the Asura_Profits workbook, per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3260.9194  # H64: 3360.1155 -> 3260.9194
$ws.Cells.Item(64, 9).Value = 3073.147  # I64: 3189.8462 -> 3073.147
$ws.Cells.Item(64, 10).Value = 3488.9285  # J64: 3530.3845 -> 3488.9285
$ws.Cells.Item(64, 11).Value = 3073.147  # K64: 3189.8462 -> 3073.147
$ws.Cells.Item(64, 12).Value = 3488.9285  # L64: 3530.3845 -> 3488.9285
$ws.Cells.Item(64, 13).Value = -2825.147  # M64: -2941.8462 -> -2825.147
$ws.Cells.Item(64, 14).Value = -3984.9285  # N64: -4026.3845 -> -3984.9285
$ws.Cells.Item(67, 8).Value = 3260.9194  # H67: 3360.1155 -> 3260.9194
$ws.Cells.Item(67, 9).Value = 3073.147  # I67: 3189.8462 -> 3073.147
$ws.Cells.Item(67, 10).Value = 3488.9285  # J67: 3530.3845 -> 3488.9285
$ws.Cells.Item(67, 11).Value = 3073.147  # K67: 3189.8462 -> 3073.147
$ws.Cells.Item(67, 12).Value = 3488.9285  # L67: 3530.3845 -> 3488.9285
$ws.Cells.Item(67, 13).Value = -2215.147  # M67: -2331.8462 -> -2215.147
$ws.Cells.Item(67, 14).Value = -5204.9285  # N67: -5246.3845 -> -5204.9285
$ws.Cells.Item(74, 8).Value = 86545.73  # H74: 79750.25 -> 86545.73
$ws.Cells.Item(74, 9).Value = 154167.17  # I74: 303334.34 -> 154167.17
$ws.Cells.Item(74, 10).Value = 5400  # J74: 5222.222 -> 5400
$ws.Cells.Item(74, 11).Value = 154167.17  # K74: 303334.34 -> 154167.17
$ws.Cells.Item(74, 12).Value = 5400  # L74: 5222.222 -> 5400
$ws.Cells.Item(74, 13).Value = -153231.17  # M74: -302398.34 -> -153231.17
$ws.Cells.Item(74, 14).Value = -7272  # N74: -7094.222 -> -7272
$ws.Cells.Item(76, 8).Value = 6411.4443  # H76: 5133.9287 -> 6411.4443
$ws.Cells.Item(76, 9).Value = 6712.875  # I76: 6571 -> 6712.875
$ws.Cells.Item(76, 10).Value = 4000  # J76: 3696.8572 -> 4000
$ws.Cells.Item(76, 11).Value = 6712.875  # K76: 6571 -> 6712.875
$ws.Cells.Item(76, 12).Value = 4000  # L76: 3696.8572 -> 4000
$ws.Cells.Item(76, 13).Value = -6397.875  # M76: -6256 -> -6397.875
$ws.Cells.Item(76, 14).Value = -4630  # N76: -4326.8572 -> -4630
$ws.Cells.Item(77, 8).Value = 86545.73  # H77: 79750.25 -> 86545.73
$ws.Cells.Item(77, 9).Value = 154167.17  # I77: 303334.34 -> 154167.17
$ws.Cells.Item(77, 10).Value = 5400  # J77: 5222.222 -> 5400
$ws.Cells.Item(77, 11).Value = 770835.8500000001  # K77: 1516671.7 -> 770835.8500000001
$ws.Cells.Item(77, 12).Value = 27000  # L77: 26111.11 -> 27000
$ws.Cells.Item(77, 13).Value = -766155.8500000001  # M77: -1511991.7 -> -766155.8500000001
$ws.Cells.Item(77, 14).Value = -36360  # N77: -35471.11 -> -36360
$ws.Cells.Item(79, 8).Value = 6411.4443  # H79: 5133.9287 -> 6411.4443
$ws.Cells.Item(79, 9).Value = 6712.875  # I79: 6571 -> 6712.875
$ws.Cells.Item(79, 10).Value = 4000  # J79: 3696.8572 -> 4000
$ws.Cells.Item(79, 11).Value = 6712.875  # K79: 6571 -> 6712.875
$ws.Cells.Item(79, 12).Value = 4000  # L79: 3696.8572 -> 4000
$ws.Cells.Item(79, 13).Value = -5620.875  # M79: -5479 -> -5620.875
$ws.Cells.Item(79, 14).Value = -6184  # N79: -5880.8572 -> -6184

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1115.7142  # H45: 1459.4667 -> 1115.7142
$ws.Cells.Item(45, 9).Value = 1122  # I45: 1124.5 -> 1122
$ws.Cells.Item(45, 10).Value = 1100  # J45: 1842.2858 -> 1100
$ws.Cells.Item(45, 11).Value = 1122  # K45: 1124.5 -> 1122
$ws.Cells.Item(45, 12).Value = 1100  # L45: 1842.2858 -> 1100
$ws.Cells.Item(45, 13).Value = -745  # M45: -747.5 -> -745
$ws.Cells.Item(45, 14).Value = -1854  # N45: -2596.2858 -> -1854
$ws.Cells.Item(74, 8).Value = 1075.4073  # H74: 1055.8572 -> 1075.4073
$ws.Cells.Item(74, 9).Value = 657.3333  # I74: 650.5263 -> 657.3333
$ws.Cells.Item(74, 11).Value = 657.3333  # K74: 650.5263 -> 657.3333
$ws.Cells.Item(74, 13).Value = 216.6667  # M74: 223.4737 -> 216.6667
$ws.Cells.Item(77, 8).Value = 1075.4073  # H77: 1055.8572 -> 1075.4073
$ws.Cells.Item(77, 9).Value = 657.3333  # I77: 650.5263 -> 657.3333
$ws.Cells.Item(77, 11).Value = 3286.6665  # K77: 3252.6315 -> 3286.6665
$ws.Cells.Item(77, 13).Value = 1081.3335  # M77: 1115.3685 -> 1081.3335
$ws.Cells.Item(96, 8).Value = 0  # H96: 37172 -> 0
$ws.Cells.Item(96, 10).Value = 0  # J96: 37172 -> 0
$ws.Cells.Item(96, 12).Value = 0  # L96: 37172 -> 0
$ws.Cells.Item(96, 14).ClearContents()  # N96 was -42664
$ws.Cells.Item(97, 8).Value = 0  # H97: 1495 -> 0
$ws.Cells.Item(97, 9).Value = 0  # I97: 1495 -> 0
$ws.Cells.Item(97, 11).Value = 0  # K97: 1495 -> 0
$ws.Cells.Item(97, 13).ClearContents()  # M97 was -999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 3350  # H94: 3233.3333 -> 3350
$ws.Cells.Item(94, 9).Value = 0  # I94: 3000 -> 0
$ws.Cells.Item(94, 11).Value = 0  # K94: 3000 -> 0
$ws.Cells.Item(94, 13).ClearContents()  # M94 was -2549
$ws.Cells.Item(107, 8).Value = 45201.5  # H107: 59333.223 -> 45201.5
$ws.Cells.Item(107, 9).Value = 53241.8  # I107: 74857 -> 53241.8
$ws.Cells.Item(107, 11).Value = 53241.8  # K107: 74857 -> 53241.8
$ws.Cells.Item(107, 13).Value = -51321.8  # M107: -72937 -> -51321.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(59, 8).Value = 44874.5  # H59: 49900 -> 44874.5
$ws.Cells.Item(59, 10).Value = 44874.5  # J59: 49900 -> 44874.5
$ws.Cells.Item(59, 12).Value = 44874.5  # L59: 49900 -> 44874.5
$ws.Cells.Item(59, 14).Value = -47164.5  # N59: -52190 -> -47164.5
$ws.Cells.Item(62, 8).Value = 38528.93  # H62: 58388.332 -> 38528.93
$ws.Cells.Item(62, 9).Value = 44283.75  # I62: 65061.875 -> 44283.75
$ws.Cells.Item(62, 10).Value = 4000  # J62: 5000 -> 4000
$ws.Cells.Item(62, 11).Value = 44283.75  # K62: 65061.875 -> 44283.75
$ws.Cells.Item(62, 12).Value = 4000  # L62: 5000 -> 4000
$ws.Cells.Item(62, 13).Value = -43659.75  # M62: -64437.875 -> -43659.75
$ws.Cells.Item(62, 14).Value = -5248  # N62: -6248 -> -5248
$ws.Cells.Item(65, 8).Value = 38528.93  # H65: 58388.332 -> 38528.93
$ws.Cells.Item(65, 9).Value = 44283.75  # I65: 65061.875 -> 44283.75
$ws.Cells.Item(65, 10).Value = 4000  # J65: 5000 -> 4000
$ws.Cells.Item(65, 11).Value = 221418.75  # K65: 325309.375 -> 221418.75
$ws.Cells.Item(65, 12).Value = 20000  # L65: 25000 -> 20000
$ws.Cells.Item(65, 13).Value = -218298.75  # M65: -322189.375 -> -218298.75
$ws.Cells.Item(65, 14).Value = -26240  # N65: -31240 -> -26240
$ws.Cells.Item(68, 8).Value = 0  # H68: 32000 -> 0
$ws.Cells.Item(68, 10).Value = 0  # J68: 32000 -> 0
$ws.Cells.Item(68, 12).Value = 0  # L68: 32000 -> 0
$ws.Cells.Item(68, 14).ClearContents()  # N68 was -33498
$ws.Cells.Item(71, 8).Value = 0  # H71: 32000 -> 0
$ws.Cells.Item(71, 10).Value = 0  # J71: 32000 -> 0
$ws.Cells.Item(71, 12).Value = 0  # L71: 96000 -> 0
$ws.Cells.Item(71, 14).ClearContents()  # N71 was -103488
$ws.Cells.Item(74, 8).Value = 0  # H74: 48000 -> 0
$ws.Cells.Item(74, 10).Value = 0  # J74: 48000 -> 0
$ws.Cells.Item(74, 12).Value = 0  # L74: 48000 -> 0
$ws.Cells.Item(74, 14).ClearContents()  # N74 was -49748
$ws.Cells.Item(77, 8).Value = 0  # H77: 48000 -> 0
$ws.Cells.Item(77, 10).Value = 0  # J77: 48000 -> 0
$ws.Cells.Item(77, 12).Value = 0  # L77: 144000 -> 0
$ws.Cells.Item(77, 14).ClearContents()  # N77 was -152736
$ws.Cells.Item(107, 8).Value = 777.5  # H107: 583.25 -> 777.5
$ws.Cells.Item(107, 9).Value = 1005.5  # I107: 590.25 -> 1005.5
$ws.Cells.Item(107, 10).Value = 549.5  # J107: 576.25 -> 549.5
$ws.Cells.Item(107, 11).Value = 1005.5  # K107: 590.25 -> 1005.5
$ws.Cells.Item(107, 12).Value = 549.5  # L107: 576.25 -> 549.5
$ws.Cells.Item(107, 13).Value = 914.5  # M107: 1329.75 -> 914.5
$ws.Cells.Item(107, 14).Value = -4389.5  # N107: -4416.25 -> -4389.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8070.3  # H70: 8637.9375 -> 8070.3
$ws.Cells.Item(70, 9).Value = 9215.846  # I70: 10734.111 -> 9215.846
$ws.Cells.Item(70, 11).Value = 9215.846  # K70: 10734.111 -> 9215.846
$ws.Cells.Item(70, 13).Value = -8945.846  # M70: -10464.111 -> -8945.846
$ws.Cells.Item(73, 8).Value = 8070.3  # H73: 8637.9375 -> 8070.3
$ws.Cells.Item(73, 9).Value = 9215.846  # I73: 10734.111 -> 9215.846
$ws.Cells.Item(73, 11).Value = 9215.846  # K73: 10734.111 -> 9215.846
$ws.Cells.Item(73, 13).Value = -8279.846  # M73: -9798.111000000001 -> -8279.846
$ws.Cells.Item(80, 8).Value = 3180  # H80: 3375 -> 3180
$ws.Cells.Item(80, 10).Value = 3450  # J80: 4500 -> 3450
$ws.Cells.Item(80, 12).Value = 3450  # L80: 4500 -> 3450
$ws.Cells.Item(80, 14).Value = -5446  # N80: -6496 -> -5446
$ws.Cells.Item(83, 8).Value = 3180  # H83: 3375 -> 3180
$ws.Cells.Item(83, 10).Value = 3450  # J83: 4500 -> 3450
$ws.Cells.Item(83, 12).Value = 17250  # L83: 22500 -> 17250
$ws.Cells.Item(83, 14).Value = -27234  # N83: -32484 -> -27234
$ws.Cells.Item(97, 8).Value = 102402.2  # H97: 52563.1 -> 102402.2
$ws.Cells.Item(97, 9).Value = 168000  # I97: 128657.75 -> 168000
$ws.Cells.Item(97, 10).Value = 4005.5  # J97: 1833.3334 -> 4005.5
$ws.Cells.Item(97, 11).Value = 168000  # K97: 128657.75 -> 168000
$ws.Cells.Item(97, 12).Value = 4005.5  # L97: 1833.3334 -> 4005.5
$ws.Cells.Item(97, 13).Value = -167504  # M97: -128161.75 -> -167504
$ws.Cells.Item(97, 14).Value = -4997.5  # N97: -2825.3334 -> -4997.5
$ws.Cells.Item(113, 8).Value = 3278.5715  # H113: 3073.7144 -> 3278.5715
$ws.Cells.Item(113, 9).Value = 1483.3334  # I113: 2331 -> 1483.3334
$ws.Cells.Item(113, 10).Value = 4625  # J113: 3816.4285 -> 4625
$ws.Cells.Item(113, 11).Value = 1483.3334  # K113: 2331 -> 1483.3334
$ws.Cells.Item(113, 12).Value = 4625  # L113: 3816.4285 -> 4625
$ws.Cells.Item(113, 13).Value = 686.6666  # M113: -161 -> 686.6666
$ws.Cells.Item(113, 14).Value = -8965  # N113: -8156.4285 -> -8965
$ws.Cells.Item(122, 8).Value = 3066.2307  # H122: 3497.3333 -> 3066.2307
$ws.Cells.Item(122, 9).Value = 2148.125  # I122: 2825 -> 2148.125
$ws.Cells.Item(122, 10).Value = 4535.2  # J122: 4035.2 -> 4535.2
$ws.Cells.Item(122, 11).Value = 6444.375  # K122: 8475 -> 6444.375
$ws.Cells.Item(122, 12).Value = 13605.6  # L122: 12105.6 -> 13605.6
$ws.Cells.Item(122, 13).Value = -3994.375  # M122: -6025 -> -3994.375
$ws.Cells.Item(122, 14).Value = -18505.6  # N122: -17005.6 -> -18505.6
$ws.Cells.Item(126, 8).Value = 2541.5652  # H126: 2987.5 -> 2541.5652
$ws.Cells.Item(126, 9).Value = 2375.2307  # I126: 2180 -> 2375.2307
$ws.Cells.Item(126, 10).Value = 2757.8  # J126: 4333.3335 -> 2757.8
$ws.Cells.Item(126, 11).Value = 7125.6921  # K126: 6540 -> 7125.6921
$ws.Cells.Item(126, 12).Value = 8273.400000000001  # L126: 13000.0005 -> 8273.400000000001
$ws.Cells.Item(126, 13).Value = -4655.6921  # M126: -4070 -> -4655.6921
$ws.Cells.Item(126, 14).Value = -13213.4  # N126: -17940.0005 -> -13213.4
$ws.Cells.Item(132, 8).Value = 3318.5  # H132: 2891.7856 -> 3318.5
$ws.Cells.Item(132, 9).Value = 2993.25  # I132: 2754.6 -> 2993.25
$ws.Cells.Item(132, 10).Value = 3535.3333  # J132: 2968 -> 3535.3333
$ws.Cells.Item(132, 11).Value = 8979.75  # K132: 8263.799999999999 -> 8979.75
$ws.Cells.Item(132, 12).Value = 10605.9999  # L132: 8904 -> 10605.9999
$ws.Cells.Item(132, 13).Value = -6449.75  # M132: -5733.799999999999 -> -6449.75
$ws.Cells.Item(132, 14).Value = -15665.9999  # N132: -13964 -> -15665.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2396.6667  # H68: 2583 -> 2396.6667
$ws.Cells.Item(68, 9).Value = 2287.5  # I68: 2536.6667 -> 2287.5
$ws.Cells.Item(68, 10).Value = 2833.3333  # J68: 3000 -> 2833.3333
$ws.Cells.Item(68, 11).Value = 2287.5  # K68: 2536.6667 -> 2287.5
$ws.Cells.Item(68, 12).Value = 2833.3333  # L68: 3000 -> 2833.3333
$ws.Cells.Item(68, 13).Value = -1538.5  # M68: -1787.6667 -> -1538.5
$ws.Cells.Item(68, 14).Value = -4331.3333  # N68: -4498 -> -4331.3333
$ws.Cells.Item(71, 8).Value = 2396.6667  # H71: 2583 -> 2396.6667
$ws.Cells.Item(71, 9).Value = 2287.5  # I71: 2536.6667 -> 2287.5
$ws.Cells.Item(71, 10).Value = 2833.3333  # J71: 3000 -> 2833.3333
$ws.Cells.Item(71, 11).Value = 11437.5  # K71: 12683.3335 -> 11437.5
$ws.Cells.Item(71, 12).Value = 14166.6665  # L71: 15000 -> 14166.6665
$ws.Cells.Item(71, 13).Value = -7693.5  # M71: -8939.333500000001 -> -7693.5
$ws.Cells.Item(71, 14).Value = -21654.6665  # N71: -22488 -> -21654.6665
$ws.Cells.Item(93, 8).Value = 2650.75  # H93: 1157.05 -> 2650.75
$ws.Cells.Item(93, 9).Value = 2301.5  # I93: 696.3125 -> 2301.5
$ws.Cells.Item(93, 11).Value = 2301.5  # K93: 696.3125 -> 2301.5
$ws.Cells.Item(93, 13).Value = -1053.5  # M93: 551.6875 -> -1053.5
$ws.Cells.Item(121, 8).Value = 49086.78  # H121: 58281 -> 49086.78
$ws.Cells.Item(121, 10).Value = 49086.78  # J121: 58281 -> 49086.78
$ws.Cells.Item(121, 12).Value = 49086.78  # L121: 58281 -> 49086.78
$ws.Cells.Item(121, 14).Value = -52580.78  # N121: -61775 -> -52580.78

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(121, 8).Value = 26289.158  # H121: 31870.857 -> 26289.158
$ws.Cells.Item(121, 10).Value = 26289.158  # J121: 31870.857 -> 26289.158
$ws.Cells.Item(121, 12).Value = 26289.158  # L121: 31870.857 -> 26289.158
$ws.Cells.Item(121, 14).Value = -29783.158  # N121: -35364.857 -> -29783.158
